$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: remove the battery-related columns (H8:L8, incl. merged J8:K8) ---
# unmerge J8:K8 first, then clear the whole H8:L8 block (values + formatting)
$ws.Range("J8:K8").UnMerge()
$ws.Range("H8:L8").Clear()

# Recreate the trailing empty (unstyled) placeholder cells T8:X8 that appear
# in the refreshed row - setting an explicit "Normal" style forces Excel to
# persist an (empty, style-less) cell record without inheriting the column's
# default wrap-text style.
$ws.Range("T8:X8").Style = "Normal"

# --- Row 15: insert a new "whl" (wheel) indicator before the existing imu/int/... list ---
# Shift I15:S15 right by one column (manually, so only row 15 is affected)
$cols = @("R","Q","P","O","N","M","L","K","J","I")
foreach ($col in $cols) {
    $src = $ws.Range($col + "15")
    $destColIndex = $src.Column() + 1
    $dest = $ws.Cells.Item(15, $destColIndex)
    $dest.Value = $src.Value()
}
$ws.Range("I15").Value = "whl"
$ws.Range("I15").WrapText = $true

# --- Sheet view: scroll so column D is the left-most visible column, and move the
#     active selection ---
$ws.Range("J36").Select()
$excel.ActiveWindow.ScrollColumn = 4
